# Add newly-sourced papers to the "source" tab, and make it the active/selected sheet
# (mirrors: "added papers to source tab")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("source")

$ws.Range("A2").Value = 'fu19'
$ws.Range("B2").Value = 'Fu_et_al_2019.pdf'
$ws.Range("C2").Value = 'Yongshuo H. Fu, Shilong Piao, Xuancheng Zhou,  Xiaojun Geng, Fanghua Hao, Yann Vitasse, Ivan A. Janssens'
$ws.Range("D2").Value = 2019
$ws.Range("E2").Value = 'Short photoperiod reduces temperature sensitivity of leaf-out in saplings of Fagus sylvatica bud not in horse chestnut'
$ws.Range("F2").Value = 'Global Change Biology'
$ws.Range("G2").Value = 25
$ws.Range("H2").Value = '1696-1703'

$ws.Range("A3").Value = 'man17'
$ws.Range("B3").Value = 'Man et al. 2017.pdf'
$ws.Range("C3").Value = 'Rongzhou Man, Pengxin Lu, Qing-Lai Dang'
$ws.Range("D3").Value = 2017
$ws.Range("E3").Value = 'Insufficient chilling effects vary among boreal tree species and chilling duration'
$ws.Range("F3").Value = 'Frontiers in Plant Science'
$ws.Range("G3").Value = 8

$ws.Range("A4").Value = 'richardson18'
$ws.Range("B4").Value = 'Richardson et al. 2018.pdf'
$ws.Range("C4").Value = 'Richardson, Andrew D. and Hufkens, Koen and Milliman, Thomas and Aubrecht, Donald M. and Furze, Morgan E. and Seyednasrollah, Bijan and Krassovski, Misha B. and Latimer, John M. and Nettles, W. Robert and Heiderman, Ryan R. and Warren, Jeffrey M. and Hanson, Paul J.'
$ws.Range("D4").Value = 2108
$ws.Range("E4").Value = 'Ecosystem warming extends vegetation activity but heightens vulnerability to cold temperatures'
$ws.Range("F4").Value = 'Nature'
$ws.Range("G4").Value = 560
$ws.Range("H4").Value = '368-371'

$ws.Range("A5").Value = 'vitra17'
$ws.Range("B5").Value = 'Vitra_et_al-2017.pdf'
$ws.Range("C5").Value = 'Amarante Vitra, Armando Lenz, Yann Vitasse'
$ws.Range("D5").Value = 2017
$ws.Range("E5").Value = 'Frost hardening and dehardening potential in temperate trees from winter to budburst'
$ws.Range("F5").Value = 'New Phytologist'
$ws.Range("G5").Value = 216
$ws.Range("H5").Value = '113-123'

$ws.Range("A6").Value = 'prevey18'
$ws.Range("B6").Value = 'prevey & Harrington 2019.pdf'
$ws.Range("C6").Value = 'Janet S. Prevey and Constance A. Harrington'
$ws.Range("D6").Value = 2018
$ws.Range("E6").Value = 'Effectiveness of winter temperatures for satisfying chilling requirements for reproductive budburst of red alder (Alnus rubra)'
$ws.Range("F6").Value = 'PeerJ'

$ws.Range("A7").Value = 'flynn18'
$ws.Range("B7").Value = 'Flynn_et_al-2018-New_Phytologist.pdf'
$ws.Range("C7").Value = 'D. F. B. Flynn & E. M. Wolkovich'
$ws.Range("D7").Value = 2018
$ws.Range("E7").Value = 'Temperature and photoperiod drive spring phenology across all species in a temperate forest community'
$ws.Range("F7").Value = 'New Phytologist'
$ws.Range("G7").Value = 219
$ws.Range("H7").Value = '1353-1362'

$ws.Range("A8").Value = 'malyshev18'
$ws.Range("B8").Value = 'Malyshev et al. 2018.pdf'
$ws.Range("C8").Value = 'Andrey V. Malyshev, Hugh A.L. Henry, Andreas Bolte, Mohammed A.S. Arfin Khan,  Juergen Kreyling'
$ws.Range("D8").Value = 2018
$ws.Range("E8").Value = 'Temporal photoperiod sensitivity and forcing requirements for budburst in temperate tree seedlings'
$ws.Range("F8").Value = 'Agriculutural and Forest Meteorology'
$ws.Range("G8").Value = 248
$ws.Range("H8").Value = '82-90'

$ws.Range("A9").Value = 'nanninga17'
$ws.Range("B9").Value = 'Nanninga et al. 2017.pdf'
$ws.Range("C9").Value = 'Claudia Nanninga, Chris R. Buyarski, Andrew M. Pretorius, Rebecca A. Montegomery'
$ws.Range("D9").Value = 2017
$ws.Range("E9").Value = 'Increased exposure to chilling advances the time to budburst in North American tree species'
$ws.Range("F9").Value = 'Tree Physiology'
$ws.Range("G9").Value = 37
$ws.Range("H9").Value = '1727-1738'

$ws.Range("A10").Value = 'anzanello16'
$ws.Range("B10").Value = 'Anzanello & Biasi 20186.pdf'
$ws.Range("C10").Value = 'Rafael Anzanello, Luiz Antonio Biasi'
$ws.Range("D10").Value = 2016
$ws.Range("E10").Value = 'Base temperature as a function of genotype: a foundation for modeling phenology of temperate fruit species'
$ws.Range("F10").Value = 'Semina: Ciencias Agrarias'
$ws.Range("G10").Value = 37
$ws.Range("H10").Value = '1811-1826'

$ws.Range("A11").Value = 'anzanello18'
$ws.Range("B11").Value = 'Anzanello et al. 2018.pdf'
$ws.Range("C11").Value = 'Rafael Anzanello, Flavio Bello Fialho, Henrique Pessoa dos Santos'
$ws.Range("D11").Value = 2018
$ws.Range("E11").Value = 'Chilling requirements and dormancy evolution in grapevine buds'
$ws.Range("F11").Value = 'Ciencia e Agrotecnologia'
$ws.Range("G11").Value = 42
$ws.Range("H11").Value = '364-371'

$ws.Range("A12").Value = 'ramos17'
$ws.Range("B12").Value = 'Ramos et al. 2018.pdf'
$ws.Range("C12").Value = 'A. Ramos, H.F. Rapoport, D. Cabello, L. Rallo'
$ws.Range("D12").Value = 2018
$ws.Range("E12").Value = 'Chilling accumulation, dormancy release temperature, and the role of leaves in olive reproductive budburst: Evaluation using shoot explants'
$ws.Range("F12").Value = 'Scientia Horticulturae'
$ws.Range("G12").Value = 231
$ws.Range("H12").Value = '241-252'

$ws.Range("A13").Value = 'fu18'
$ws.Range("B13").Value = 'Fu_et_al-2018-Global_Change_Biology.pdf'
$ws.Range("C13").Value = 'Yongshuo H. Fu, Shilong Piao, Nicolas Delpierre, Fanghua Hao, Heikki H€anninen, Yongjie Liu, Wenchao Sun, Ivan A. Janssens, Matteo Campioli'
$ws.Range("D13").Value = 2018
$ws.Range("E13").Value = 'Larger temperature response of autumn leaf senescence than spring leaf-out phenology'
$ws.Range("F13").Value = 'Global Change Biology'
$ws.Range("G13").Value = 24
$ws.Range("H13").Value = '2159-2168'

# Cell F10 (journal "Semina: Ciencias Agrarias") is center-aligned
$ws.Range("F10").HorizontalAlignment = -4108

# Cell H6 carries a leftover date-ish number format even though it has no value
$ws.Range("H6").NumberFormat = "d-mmm"

# Column B was widened to fit the new filenames
$ws.Columns.Item(2).AutoFit()

# Make "source" the active/selected sheet (matches tabSelected/activeTab in the diff)
$ws.Activate()
$ws.Range("B5").Select()
